$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.297.23"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").Value = "1.730.02"
$ws.Range("E3").Value = "  -3.82%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.04"
$ws.Range("E5").Value = "  -4.54%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4226"
$ws.Range("E7").Value = "  -10.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3576"
$ws.Range("E8").Value = "  -3.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.86"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07403"
$ws.Range("E10").Value = "  -3.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.104"
$ws.Range("E11").Value = "  -3.83%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.39"
$ws.Range("E13").Value = "  -5.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.055"
$ws.Range("E14").Value = "  -4.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.101"
$ws.Range("E15").Value = "  -3.95%  "
$ws.Range("D16").Value = "1.729.80"
$ws.Range("E16").Value = "  -3.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001056"
$ws.Range("E17").Value = "  -3.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.65"
$ws.Range("E18").Value = "  +4.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05942"
$ws.Range("E19").Value = "  -11.92%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.67"
$ws.Range("E21").Value = "  -4.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.062"
$ws.Range("E22").Value = "  -5.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5239"
$ws.Range("E23").Value = "  -5.76%  "
$ws.Range("D24").Value = "27.332.79"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("E25").Value = "  -5.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.392"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.05"
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.343"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "148.41"
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("D30").Value = "1.925.76"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.213"
$ws.Range("E31").Value = "  -3.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "125.69"
$ws.Range("E32").Value = "  -6.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09077"
$ws.Range("E33").Value = "  -6.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.588"
$ws.Range("E34").Value = "  -5.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.626"
$ws.Range("E35").Value = "  -10.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.69"
$ws.Range("E36").Value = "  +4.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2153"
$ws.Range("E37").Value = "  -3.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06085"
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.037"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02242"
$ws.Range("E40").Value = "  -5.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6337"
$ws.Range("E41").Value = "  -5.59%  "
$ws.Range("E42").Value = "  -4.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.412"
$ws.Range("E44").Value = "  -6.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.882"
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("E46").Value = "  -4.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.726"
$ws.Range("E47").Value = "  -3.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5796"
$ws.Range("E48").Value = "  -5.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.53"
$ws.Range("E49").Value = "  -4.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.937"
$ws.Range("E50").Value = "  -6.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06811"
$ws.Range("E51").Value = "  -4.55%  "
